$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsForecast.Range("D13").Value = 157
$wsForecast.Range("D14").Value = 152
$wsForecast.Range("D15").Value = 158
$wsForecast.Range("D16").Value = 140
$wsForecast.Range("D17").Value = 103

# --- Sheet: Summary ---
# These cells hold numeric-looking values stored as text (inline strings)
# in the source workbook, so force text entry then strip the quote-prefix
# style that Excel auto-applies, keeping the cell on the default style.
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "2899"
$wsSummary.Range("B9").ClearFormats()

$wsSummary.Range("B14").NumberFormat = "@"
$wsSummary.Range("B14").Value = "103"
$wsSummary.Range("B14").ClearFormats()
